$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '98.907.31'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '3.323.36'
$ws.Range("E3").Value = '  -1.22%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.11'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.35%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '632.84'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.62%  '
$ws.Range("E7").Value = '  +18.16%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.411'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +6.03%  '
$ws.Range("E9").Value = '  -0.04%  '
$ws.Range("E10").Value = '  +22.72%  '
$ws.Range("D11").Value = '3.319.47'
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("E12").Value = '  +3.00%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '43.11'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +20.40%  '
$ws.Range("D14").Value = '98.649.70'
$ws.Range("E14").Value = '  +0.55%  '
$ws.Range("E15").Value = '  +2.41%  '
$ws.Range("D16").Value = '3.960.26'
$ws.Range("E16").Value = '  -0.92%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.44'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.36%  '
$ws.Range("D18").Value = '3.322.69'
$ws.Range("E18").Value = '  -1.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.54'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +10.82%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.56'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -2.15%  '
$ws.Range("E21").Value = '  +11.94%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '488.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.45%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.76'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.24%  '
$ws.Range("E24").Value = '  -1.04%  '
$ws.Range("E25").Value = '  +7.52%  '
$ws.Range("E26").Value = '  +31.26%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '89.99'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.41'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("D29").Value = '3.503.94'
$ws.Range("E29").Value = '  -1.13%  '
$ws.Range("E30").Value = '  +16.62%  '
$ws.Range("B31").Value = 'Dai'
$ws.Range("C31").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.999'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.06%  '
$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +21.56%  '
$ws.Range("E33").Value = '  +2.75%  '
$ws.Range("E34").Value = '  +0.00%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '28.37'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.63%  '
$ws.Range("E36").Value = '  +8.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '7.55'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("E38").Value = '  -0.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.39%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '510.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '24.77'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("E42").Value = '  +2.65%  '
$ws.Range("E43").Value = '  -0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.809'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +3.57%  '
$ws.Range("E45").Value = '  +0.06%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.20'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.04'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.79%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '160.81'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.30%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.48'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +8.23%  '
$ws.Range("B50").Value = 'Filecoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '4.91'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +8.21%  '
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +16.21%  '
